$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the numeric section identifiers in column A with the
# corresponding "Стол N" (Table N) text labels.
$ws.Range("A2:A6").Value = "Стол 1  "
$ws.Range("A7:A11").Value = "Стол 2"
$ws.Range("A12:A15").Value = "Стол 3"
$ws.Range("A16:A20").Value = "Стол 4"
$ws.Range("A21:A24").Value = "Стол 5"

# Update the active selection to match the saved workbook state.
$ws.Range("A22").Select()
